# Sprint 4 Backlog Burndown - Week 1 progress update

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fixes section: record Week 1 remaining amount for "Fixes" task (row 3)
$ws.Range("E3").Value = 0

# Tickets section: assign owners and record estimate / week 1 remaining
# Create Ticket -> Jacob, estimate 4, week1 remaining 4
$ws.Range("C9").Value = "Jacob"
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = 4

# Update Ticket Stage -> Jabesi
$ws.Range("C10").Value = "Jabesi"

# Rename "Assign Ticket" task to "Assign Ticket to User", assign to Jacob,
# estimate 4, week1 remaining 4
$ws.Range("B11").Value = "Assign Ticket to User"
$ws.Range("C11").Value = "Jacob"
$ws.Range("D11").Value = 4
$ws.Range("E11").Value = 4

# Update selection to reflect where the user left off editing
$ws.Range("H10").Select()

$wb.Save()
